$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record is inserted at row 22, shifting all existing
# records (rows 22:173) down by one row (the old row 173 becomes row 174).
$ws.Rows(22).Insert()

# Populate the newly inserted row 22 with the new record's data.
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value = 44462
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 100114014
$ws.Range("G22").Value = "Betarraga"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = 650
$ws.Range("N22").Value = "$/paquete 5 unidades"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 130
$ws.Range("Q22").Value = 5
$ws.Range("R22").Value = "Hortaliza"
